$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 9.712432666666666
$ws.Range("H2").Value = 29.137298
$ws.Range("I2").Value = 0.4639063029983291
$ws.Range("J2").Value = 0.463906302998329
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.709791333333333
$ws.Range("N2").Value = 11.129374
$ws.Range("O2").Value = 0.4283284425582907
$ws.Range("P2").Value = 0.4283284425582907
$ws.Range("Q2").Value = 36.03109853238355
$ws.Range("R2").Value = 324.279886791452
$ws.Range("S2").Value = 0.1987042642562488
$ws.Range("T2").Value = 0.1987042642562488

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 9.712432666666666
$ws.Range("H3").Value = 29.137298
$ws.Range("I3").Value = 0.4639063029983291
$ws.Range("J3").Value = 0.463906302998329
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.599001333333334
$ws.Range("N3").Value = 10.797004
$ws.Range("O3").Value = 0.4155367505499981
$ws.Range("P3").Value = 0.4155367505499982
$ws.Range("Q3").Value = 34.95505811724356
$ws.Range("R3").Value = 314.595523055192
$ws.Range("S3").Value = 0.1927701177075885
$ws.Range("T3").Value = 0.1927701177075885

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 9.712432666666666
$ws.Range("H4").Value = 29.137298
$ws.Range("I4").Value = 0.4639063029983291
$ws.Range("J4").Value = 0.463906302998329
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.352297666666667
$ws.Range("N4").Value = 4.056893
$ws.Range("O4").Value = 0.1561348068917112
$ws.Range("P4").Value = 0.1561348068917112
$ws.Range("Q4").Value = 13.13410003279044
$ws.Range("R4").Value = 118.206900295114
$ws.Range("S4").Value = 0.07243192103449178
$ws.Range("T4").Value = 0.07243192103449177

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 9.124904999999998
$ws.Range("H5").Value = 27.37471499999999
$ws.Range("I5").Value = 0.4358435305594535
$ws.Range("J5").Value = 0.4358435305594534
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.709791333333333
$ws.Range("N5").Value = 11.129374
$ws.Range("O5").Value = 0.4283284425582907
$ws.Range("P5").Value = 0.4283284425582907
$ws.Range("Q5").Value = 33.85149348648999
$ws.Range("R5").Value = 304.66344137841
$ws.Range("S5").Value = 0.1866841806436375
$ws.Range("T5").Value = 0.1866841806436375

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 9.124904999999998
$ws.Range("H6").Value = 27.37471499999999
$ws.Range("I6").Value = 0.4358435305594535
$ws.Range("J6").Value = 0.4358435305594534
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.599001333333334
$ws.Range("N6").Value = 10.797004
$ws.Range("O6").Value = 0.4155367505499981
$ws.Range("P6").Value = 0.4155367505499982
$ws.Range("Q6").Value = 32.84054526153999
$ws.Range("R6").Value = 295.56490735386
$ws.Range("S6").Value = 0.1811090044369141
$ws.Range("T6").Value = 0.1811090044369141

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 9.124904999999998
$ws.Range("H7").Value = 27.37471499999999
$ws.Range("I7").Value = 0.4358435305594535
$ws.Range("J7").Value = 0.4358435305594534
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.352297666666667
$ws.Range("N7").Value = 4.056893
$ws.Range("O7").Value = 0.1561348068917112
$ws.Range("P7").Value = 0.1561348068917112
$ws.Range("Q7").Value = 12.339587740055
$ws.Range("R7").Value = 111.056289660495
$ws.Range("S7").Value = 0.06805034547890189
$ws.Range("T7").Value = 0.06805034547890189

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.318184
$ws.Range("H8").Value = 0.9545520000000001
$ws.Range("I8").Value = 0.01519779525677573
$ws.Range("J8").Value = 0.01519779525677573
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.709791333333333
$ws.Range("N8").Value = 11.129374
$ws.Range("O8").Value = 0.4283284425582907
$ws.Range("P8").Value = 0.4283284425582907
$ws.Range("Q8").Value = 1.180396245605333
$ws.Range("R8").Value = 10.623566210448
$ws.Range("S8").Value = 0.006509647972654528
$ws.Range("T8").Value = 0.006509647972654528

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.318184
$ws.Range("H9").Value = 0.9545520000000001
$ws.Range("I9").Value = 0.01519779525677573
$ws.Range("J9").Value = 0.01519779525677573
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.599001333333334
$ws.Range("N9").Value = 10.797004
$ws.Range("O9").Value = 0.4155367505499981
$ws.Range("P9").Value = 0.4155367505499982
$ws.Range("Q9").Value = 1.145144640245334
$ws.Range("R9").Value = 10.306301762208
$ws.Range("S9").Value = 0.006315242456524763
$ws.Range("T9").Value = 0.006315242456524763

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.318184
$ws.Range("H10").Value = 0.9545520000000001
$ws.Range("I10").Value = 0.01519779525677573
$ws.Range("J10").Value = 0.01519779525677573
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.352297666666667
$ws.Range("N10").Value = 4.056893
$ws.Range("O10").Value = 0.1561348068917112
$ws.Range("P10").Value = 0.1561348068917112
$ws.Range("Q10").Value = 0.4302794807706667
$ws.Range("R10").Value = 3.872515326936
$ws.Range("S10").Value = 0.002372904827596444
$ws.Range("T10").Value = 0.002372904827596444

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.353022
$ws.Range("H11").Value = 1.059066
$ws.Range("I11").Value = 0.01686180347577968
$ws.Range("J11").Value = 0.01686180347577968
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 3.709791333333333
$ws.Range("N11").Value = 11.129374
$ws.Range("O11").Value = 0.4283284425582907
$ws.Range("P11").Value = 0.4283284425582907
$ws.Range("Q11").Value = 1.309637956076
$ws.Range("R11").Value = 11.786741604684
$ws.Range("S11").Value = 0.007222390021504685
$ws.Range("T11").Value = 0.007222390021504684

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.353022
$ws.Range("H12").Value = 1.059066
$ws.Range("I12").Value = 0.01686180347577968
$ws.Range("J12").Value = 0.01686180347577968
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 3.599001333333334
$ws.Range("N12").Value = 10.797004
$ws.Range("O12").Value = 0.4155367505499981
$ws.Range("P12").Value = 0.4155367505499982
$ws.Range("Q12").Value = 1.270526648696
$ws.Range("R12").Value = 11.434739838264
$ws.Range("S12").Value = 0.007006699024738154
$ws.Range("T12").Value = 0.007006699024738153

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.353022
$ws.Range("H13").Value = 1.059066
$ws.Range("I13").Value = 0.01686180347577968
$ws.Range("J13").Value = 0.01686180347577968
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.352297666666667
$ws.Range("N13").Value = 4.056893
$ws.Range("O13").Value = 0.1561348068917112
$ws.Range("P13").Value = 0.1561348068917112
$ws.Range("Q13").Value = 0.477390826882
$ws.Range("R13").Value = 4.296517441938
$ws.Range("S13").Value = 0.002632714429536846
$ws.Range("T13").Value = 0.002632714429536845

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 1.427651
$ws.Range("H14").Value = 4.282953
$ws.Range("I14").Value = 0.06819056770966213
$ws.Range("J14").Value = 0.06819056770966211
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 3.709791333333333
$ws.Range("N14").Value = 11.129374
$ws.Range("O14").Value = 0.4283284425582907
$ws.Range("P14").Value = 0.4283284425582907
$ws.Range("Q14").Value = 5.296287306824667
$ws.Range("R14").Value = 47.666585761422
$ws.Range("S14").Value = 0.02920795966424524
$ws.Range("T14").Value = 0.02920795966424524

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 1.427651
$ws.Range("H15").Value = 4.282953
$ws.Range("I15").Value = 0.06819056770966213
$ws.Range("J15").Value = 0.06819056770966211
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 3.599001333333334
$ws.Range("N15").Value = 10.797004
$ws.Range("O15").Value = 0.4155367505499981
$ws.Range("P15").Value = 0.4155367505499982
$ws.Range("Q15").Value = 5.138117852534667
$ws.Range("R15").Value = 46.24306067281201
$ws.Range("S15").Value = 0.02833568692423263
$ws.Range("T15").Value = 0.02833568692423263

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 1.427651
$ws.Range("H16").Value = 4.282953
$ws.Range("I16").Value = 0.06819056770966213
$ws.Range("J16").Value = 0.06819056770966211
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1.352297666666667
$ws.Range("N16").Value = 4.056893
$ws.Range("O16").Value = 0.1561348068917112
$ws.Range("P16").Value = 0.1561348068917112
$ws.Range("Q16").Value = 1.930609116114333
$ws.Range("R16").Value = 17.375482045029
$ws.Range("S16").Value = 0.01064692112118425
$ws.Range("T16").Value = 0.01064692112118425

